$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.671.95"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "3.564.92"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'607.09"
$ws.Range("E5").Value = "  +4.05%  "
$ws.Range("D6").Value = "'173.82"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "3.561.74"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "'0.196"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").Value = "'7.36"
$ws.Range("E11").Value = "  +8.67%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'46.80"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'0.0000277"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "4.142.94"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").Value = "'616.27"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").Value = "3.564.20"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "70.767.65"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "'0.888"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'9.40"
$ws.Range("E23").Value = "  -16.81%  "
$ws.Range("D24").Value = "'16.02"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "'97.39"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "'3.81"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "'9.11"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D31").Value = "'8.48"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("D35").Value = "'611.67"
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("D36").Value = "'3.73"
$ws.Range("E36").Value = "  +5.41%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'0.0481"
$ws.Range("E39").Value = "  +5.20%  "
$ws.Range("D40").Value = "'57.33"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("D43").Value = "3.384.36"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("E45").Value = "  +7.84%  "
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0₃0707"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'32.98"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "'132.94"
$ws.Range("E50").Value = "  -0.09%  "
